# major update: add column to database
# Adds a new "pengguna" (user) column at the front of the traffic-log
# table and replaces the data rows with the latest 3 log entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so stale rows (5:7 in the old layout) disappear.
$ws.Cells.Clear()

# Header row (A1:H1) — "pengguna" is the newly added column.
$ws.Range("A1").Value = "pengguna"
$ws.Range("B1").Value = "Interval"
$ws.Range("C1").Value = "Durasi"
$ws.Range("D1").Value = "SM"
$ws.Range("E1").Value = "MP"
$ws.Range("F1").Value = "KS"
$ws.Range("G1").Value = "BB"
$ws.Range("H1").Value = "TB"

# Row 2 — no user recorded for this entry (column A left blank).
$ws.Range("B2").Value = "3/13/2024 , 12:12:09 - 12:12:14"
$ws.Range("C2").Value = "0 Menit 5 Detik"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0

# Row 3 — user "wiki".
$ws.Range("A3").Value = "wiki"
$ws.Range("B3").Value = "3/13/2024 , 12:55:30 - 12:55:33"
$ws.Range("C3").Value = "0 Menit 3 Detik"
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 11
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 31

# Row 4 — user "mikhael".
$ws.Range("A4").Value = "mikhael"
$ws.Range("B4").Value = "3/13/2024 , 12:56:57 - 12:57:05"
$ws.Range("C4").Value = "0 Menit 8 Detik"
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 1
